# moved local logs to repo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the weekly summary header info
$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 3

# Stage / Task names first (column A, then column B) to mirror authoring order
$ws.Range("A3").Value = "Project Planning"
$ws.Range("A4").Value = "Analysis/requirements Elicitation"
$ws.Range("B3").Value = "Finalise Project Plan"
$ws.Range("B4").Value = "Client Meeting, Gather and analyse requirements"

# Hours columns
$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 18
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2

# Total hours spent this week
$ws.Range("B14").Value = 60

# Widen column A to fit the new text (engine's ColumnWidth setter rounds to the
# stored character-width grid; 31.15 converges to a stored width of 32)
$ws.Columns.Item(1).ColumnWidth = 31.15

# Update the active selection to match the saved view
$ws.Range("B13").Select()
